$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Deje" paragraph - merge the run that held the word "si" (which
# was wrapped in w:proofErr gramStart/gramEnd markers) back into the
# surrounding text, producing one continuous run/sentence.
# Re-typing the exact same text over the whole affected span collapses the
# three runs (and drops the proofErr grammar markers) into a single run.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(
    "myslí si že je skutečně mrtvá, u její hrobky probodne Parise, sám vypije jed a zemře, Julie se probudí a probodne se dýkou.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "myslí si že je skutečně mrtvá, u její hrobky probodne Parise, sám vypije jed a zemře, Julie se probudí a probodne se dýkou.",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Giovanny Boccacio" entry under "Italska renesancni literatura"
# becomes a proper Nadpis4 (heading 4) entry, named correctly ("Giovanni"),
# with "Boccacio" underlined (in addition to already being bold), and the
# whole entry (plus the following "Dilo: ..." bullet) wrapped by a bookmark.
# ---------------------------------------------------------------------------

# Locate the paragraph that currently reads "Giovanny Boccacio".
$findRng = $d.Content
$findRng.Find.Execute("Giovanny", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$giovanniStart = $findRng.Start
$giovanniPara = $findRng.Paragraphs(1)

# Switch it from the bulleted "List Paragraph" style to "Nadpis4" - this also
# drops the direct numbering/spacing formatting that belonged to the list.
$giovanniPara.Style = "Nadpis4"

# Fix the typo "Giovanny" -> "Giovanni" while re-creating the 3 separate runs
# ("Giovann" / "i" / " ") that the target markup expects.
$wholeWord = $d.Range($giovanniStart, $giovanniStart + 8)
$wholeWord.Text = "Giovanni"
$iChar = $d.Range($giovanniStart + 7, $giovanniStart + 8)
$iChar.Font.Bold = 1
$iChar.Font.Bold = 0

# Underline "Boccacio" (it is already bold).
$boccacioRng = $d.Content
$boccacioRng.Find.Execute("Boccacio", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boccacioRng.Font.Underline = 1

# Wrap the heading paragraph together with the following "Dilo: Dekameron..."
# paragraph in a bookmark, matching the source document's existing bookmarks.
$diloRng = $d.Content
$diloRng.Find.Execute("Dílo: Dekameron", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$diloPara = $diloRng.Paragraphs(1)
$bookmarkRange = $d.Range($giovanniStart, $diloPara.Range.End)
$d.Bookmarks.Add("_Hlk130320128", $bookmarkRange) | Out-Null

Write-Output "done"
